$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.3949551318557
$ws.Range("C2").Value = 8.772861809360174
$ws.Range("D2").Value = 3.838276894093157
$ws.Range("F2").Value = 19.49238989544178
$ws.Range("G2").Value = 3.594448496529577
$ws.Range("N2").Value = 15.81695941436856
$ws.Range("O2").Value = 17.34286194081059
$ws.Range("B3").Value = 11.78358033903897
$ws.Range("C3").Value = 8.369125340180023
$ws.Range("D3").Value = 3.768808920947571
$ws.Range("F3").Value = 19.47647288460611
$ws.Range("G3").Value = 3.596321630334197
$ws.Range("N3").Value = 15.84310144985067
$ws.Range("O3").Value = 17.39577011495001
$ws.Range("B4").Value = 11.39219079140143
$ws.Range("C4").Value = 8.109650073007284
$ws.Range("D4").Value = 3.725070534068491
$ws.Range("F4").Value = 19.47405995146817
$ws.Range("G4").Value = 3.597533071982244
$ws.Range("N4").Value = 15.86099454005563
$ws.Range("O4").Value = 17.43381126767999
$ws.Range("B5").Value = 11.22885844739228
$ws.Range("C5").Value = 8.00109133336843
$ws.Range("D5").Value = 3.706987782928316
$ws.Range("F5").Value = 19.47492676981049
$ws.Range("G5").Value = 3.598042213923302
$ws.Range("N5").Value = 15.86874990615548
$ws.Range("O5").Value = 17.45070385563394
$ws.Range("B6").Value = 11.20151148818259
$ws.Range("C6").Value = 7.982897851722111
$ws.Range("D6").Value = 3.703969926499477
$ws.Range("F6").Value = 19.47518238683791
$ws.Range("G6").Value = 3.598127692287191
$ws.Range("N6").Value = 15.87006570936817
$ws.Range("O6").Value = 17.45359264426064
$ws.Range("B7").Value = 11.3900033051759
$ws.Range("C7").Value = 8.108197300237235
$ws.Range("D7").Value = 3.724827693845392
$ws.Range("F7").Value = 19.47406415307703
$ws.Range("G7").Value = 3.597539875748785
$ws.Range("N7").Value = 15.86109725292442
$ws.Range("O7").Value = 17.43403346572819
$ws.Range("B8").Value = 12.18757691061988
$ws.Range("C8").Value = 8.636107064793427
$ws.Range("D8").Value = 3.814557557619335
$ws.Range("F8").Value = 19.4853743197102
$ws.Range("G8").Value = 3.595081652819814
$ws.Range("N8").Value = 15.82559132518008
$ws.Range("O8").Value = 17.35994794106133
$ws.Range("B9").Value = 13.61791228746761
$ws.Range("C9").Value = 9.576152435320511
$ws.Range("D9").Value = 3.981337072348366
$ws.Range("F9").Value = 19.56588824727848
$ws.Range("G9").Value = 3.590745531114867
$ws.Range("N9").Value = 15.77054930045677
$ws.Range("O9").Value = 17.25901291558897
$ws.Range("B10").Value = 14.57999082515763
$ws.Range("C10").Value = 10.2051948350863
$ws.Range("D10").Value = 4.0975521704718
$ws.Range("F10").Value = 19.66037860331316
$ws.Range("G10").Value = 3.587852060946457
$ws.Range("N10").Value = 15.73896304994329
$ws.Range("O10").Value = 17.21222636261942
$ws.Range("B11").Value = 14.99726177108974
$ws.Range("C11").Value = 10.47743061470788
$ws.Range("D11").Value = 4.148910342505072
$ws.Range("F11").Value = 19.71094499626637
$ws.Range("G11").Value = 3.586598562668104
$ws.Range("N11").Value = 15.72650765016268
$ws.Range("O11").Value = 17.19694619571258
$ws.Range("B12").Value = 15.15226830916677
$ws.Range("C12").Value = 10.57848205739102
$ws.Range("D12").Value = 4.168130654568144
$ws.Range("F12").Value = 19.73117234163529
$ws.Range("G12").Value = 3.586132870579031
$ws.Range("N12").Value = 15.72206552327475
$ws.Range("O12").Value = 17.1920273408828
$ws.Range("B13").Value = 15.11901955285182
$ws.Range("C13").Value = 10.55680999133462
$ws.Range("D13").Value = 4.164001520148441
$ws.Range("F13").Value = 19.72676823739538
$ws.Range("G13").Value = 3.586232766984423
$ws.Range("N13").Value = 15.72301001892141
$ws.Range("O13").Value = 17.19304806632179
$ws.Range("B14").Value = 15.01007484942866
$ws.Range("C14").Value = 10.4857852158623
$ws.Range("D14").Value = 4.150496237196406
$ws.Range("F14").Value = 19.71258756088111
$ws.Range("G14").Value = 3.586560070166559
$ws.Range("N14").Value = 15.72613669700042
$ws.Range("O14").Value = 17.19652411273313
$ws.Range("B15").Value = 14.94294978846035
$ws.Range("C15").Value = 10.44201406179743
$ws.Range("D15").Value = 4.142193866865456
$ws.Range("F15").Value = 19.70404162527559
$ws.Range("G15").Value = 3.586761721051535
$ws.Range("N15").Value = 15.72808760040753
$ws.Range("O15").Value = 17.19876636698163
$ws.Range("B16").Value = 14.55230439261082
$ws.Range("C16").Value = 10.18712019733335
$ws.Range("D16").Value = 4.094164481277268
$ws.Range("F16").Value = 19.65722562075569
$ws.Range("G16").Value = 3.587935239070773
$ws.Range("N16").Value = 15.73981548383685
$ws.Range("O16").Value = 17.21334612725255
$ws.Range("B17").Value = 14.30737838563315
$ws.Range("C17").Value = 10.0271576794699
$ws.Range("D17").Value = 4.064305576740109
$ws.Range("F17").Value = 19.63044009149671
$ws.Range("G17").Value = 3.588671196422972
$ws.Range("N17").Value = 15.74749974060709
$ws.Range("O17").Value = 17.22383106460238
$ws.Range("B18").Value = 14.16458828157495
$ws.Range("C18").Value = 9.933843090191512
$ws.Range("D18").Value = 4.046990133647161
$ws.Range("F18").Value = 19.61574825633779
$ws.Range("G18").Value = 3.589100409091545
$ws.Range("N18").Value = 15.752099656362
$ws.Range("O18").Value = 17.23042656908668
$ws.Range("B19").Value = 14.11591555817196
$ws.Range("C19").Value = 9.902024895683006
$ws.Range("D19").Value = 4.041103483887643
$ws.Range("F19").Value = 19.61089688041026
$ws.Range("G19").Value = 3.589246749420262
$ws.Range("N19").Value = 15.75368806728807
$ws.Range("O19").Value = 17.23275657002865
$ws.Range("B20").Value = 14.33364996587564
$ws.Range("C20").Value = 10.04432166372719
$ws.Range("D20").Value = 4.067498823220273
$ws.Range("F20").Value = 19.63321758521292
$ws.Range("G20").Value = 3.58859224121699
$ws.Range("N20").Value = 15.74666309992359
$ws.Range("O20").Value = 17.22265643368851
$ws.Range("B21").Value = 15.042156622661
$ws.Range("C21").Value = 10.50670249382945
$ws.Range("D21").Value = 4.154469335601774
$ws.Range("F21").Value = 19.71672358782583
$ws.Range("G21").Value = 3.586463689905982
$ws.Range("N21").Value = 15.72521087248786
$ws.Range("O21").Value = 17.19547954227737
$ws.Range("B22").Value = 15.48766672345959
$ws.Range("C22").Value = 10.7970016151259
$ws.Range("D22").Value = 4.20997654801986
$ws.Range("F22").Value = 19.77758180851043
$ws.Range("G22").Value = 3.585124883797746
$ws.Range("N22").Value = 15.71279005336308
$ws.Range("O22").Value = 17.18277530531245
$ws.Range("B23").Value = 15.25151545111302
$ws.Range("C23").Value = 10.64316221842164
$ws.Range("D23").Value = 4.180476720336497
$ws.Range("F23").Value = 19.74453009618871
$ws.Range("G23").Value = 3.585834656401109
$ws.Range("N23").Value = 15.71927315622194
$ws.Range("O23").Value = 17.18909182886657
$ws.Range("B24").Value = 14.32177874247272
$ws.Range("C24").Value = 10.03656602859437
$ws.Range("D24").Value = 4.066055620250189
$ws.Range("F24").Value = 19.6319596755848
$ws.Range("G24").Value = 3.588627917858511
$ws.Range("N24").Value = 15.74704077778405
$ws.Range("O24").Value = 17.22318571667961
$ws.Range("B25").Value = 13.24611945128668
$ws.Range("C25").Value = 9.332446770950185
$ws.Range("D25").Value = 3.937278672855995
$ws.Range("F25").Value = 19.53787680651016
$ws.Range("G25").Value = 3.591867020534107
$ws.Range("N25").Value = 15.78388209241668
$ws.Range("O25").Value = 17.2815342395526
